$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.984.78"
$ws.Range("E2").Value = "  -0.02%  "

$ws.Range("D3").Value = "2.452.39"
$ws.Range("E3").Value = "  +0.14%  "

$ws.Range("E4").Value = "  -0.07%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "510.53"
$ws.Range("E5").Value = "  -2.56%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "133.72"
$ws.Range("E6").Value = "  +3.09%  "

$ws.Range("E7").Value = "  -0.04%  "

$ws.Range("E8").Value = "  -1.49%  "

$ws.Range("D9").Value = "2.452.56"
$ws.Range("E9").Value = "  -0.11%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0981"
$ws.Range("E10").Value = "  +0.45%  "

$ws.Range("E11").Value = "  -0.71%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.324"
$ws.Range("E12").Value = "  +0.46%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.61"
$ws.Range("E13").Value = "  -7.08%  "

$ws.Range("D14").Value = "2.888.92"
$ws.Range("E14").Value = "  +0.03%  "

$ws.Range("D15").Value = "57.858.67"
$ws.Range("E15").Value = "  -0.13%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "21.93"
$ws.Range("E16").Value = "  +1.72%  "

$ws.Range("E17").Value = "  +1.37%  "

$ws.Range("D18").Value = "2.465.18"
$ws.Range("E18").Value = "  +0.29%  "

$ws.Range("E19").Value = "  -0.48%  "

$ws.Range("E20").Value = "  +0.83%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "314.75"
$ws.Range("E21").Value = "  +0.23%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.44"
$ws.Range("E22").Value = "  +4.62%  "

$ws.Range("E23").Value = "  +0.02%  "

$ws.Range("E24").Value = "  -2.23%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "65.41"
$ws.Range("E25").Value = "  +0.37%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("E26").Value = "  -0.16%  "

$ws.Range("E27").Value = "  -0.79%  "

$ws.Range("E28").Value = "  -6.00%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.57"
$ws.Range("E29").Value = "  +4.38%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "172.47"
$ws.Range("E30").Value = "  -1.35%  "

$ws.Range("D31").Value = "0.0₃0733"
$ws.Range("E31").Value = "  -0.47%  "

$ws.Range("E32").Value = "  +0.13%  "

$ws.Range("E33").Value = "  +0.17%  "

$ws.Range("E34").Value = "  +0.03%  "

$ws.Range("E35").Value = "  +0.06%  "

$ws.Range("E36").Value = "  +0.13%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "18.09"
$ws.Range("E37").Value = "  +1.22%  "

$ws.Range("E38").Value = "  +5.39%  "

$ws.Range("E39").Value = "  +1.96%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "36.73"
$ws.Range("E40").Value = "  +1.32%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.46"
$ws.Range("E41").Value = "  +1.10%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.805"
$ws.Range("E42").Value = "  -0.33%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "136.85"
$ws.Range("E43").Value = "  +8.06%  "

$ws.Range("E44").Value = "  +0.38%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.90"
$ws.Range("E45").Value = "  +2.26%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.577"
$ws.Range("E46").Value = "  -1.24%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "256.52"
$ws.Range("E47").Value = "  -1.19%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0920"
$ws.Range("E48").Value = "  -0.30%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0492"
$ws.Range("E49").Value = "  -0.06%  "

$ws.Range("E50").Value = "  +1.75%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "17.17"
$ws.Range("E51").Value = "  +0.82%  "
